$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 13.889835214901
$ws.Cells.Item(2, 3).Value = 11.56967659086022
$ws.Cells.Item(2, 4).Value = 4.007240567040604
$ws.Cells.Item(2, 6).Value = 20.20909231808771
$ws.Cells.Item(2, 7).Value = 21.91095700889505
$ws.Cells.Item(2, 8).Value = 12.42426831528762
$ws.Cells.Item(2, 9).Value = 19.2168517511913
$ws.Cells.Item(2, 12).Value = 10.61908329509709
$ws.Cells.Item(2, 14).Value = 16.46556401207404
$ws.Cells.Item(2, 15).Value = 18.0108287304179
$ws.Cells.Item(3, 2).Value = 13.31967578700308
$ws.Cells.Item(3, 3).Value = 11.46732606047098
$ws.Cells.Item(3, 4).Value = 3.934638899796703
$ws.Cells.Item(3, 6).Value = 20.17811085246719
$ws.Cells.Item(3, 7).Value = 21.82759489501826
$ws.Cells.Item(3, 8).Value = 12.45914635260138
$ws.Cells.Item(3, 9).Value = 19.30845552554066
$ws.Cells.Item(3, 12).Value = 10.58939611999932
$ws.Cells.Item(3, 14).Value = 16.49807403837879
$ws.Cells.Item(3, 15).Value = 18.04553542684458
$ws.Cells.Item(4, 2).Value = 12.95776369799554
$ws.Cells.Item(4, 3).Value = 11.4039979103941
$ws.Cells.Item(4, 4).Value = 3.888789959621477
$ws.Cells.Item(4, 6).Value = 20.16537637974609
$ws.Cells.Item(4, 7).Value = 21.78524422229777
$ws.Cells.Item(4, 8).Value = 12.48283609312157
$ws.Cells.Item(4, 9).Value = 19.36903535652761
$ws.Cells.Item(4, 12).Value = 10.57328946317169
$ws.Cells.Item(4, 14).Value = 16.51973082318207
$ws.Cells.Item(4, 15).Value = 18.07152697674239
$ws.Cells.Item(5, 2).Value = 12.80751241374946
$ws.Cells.Item(5, 3).Value = 11.37808564609683
$ws.Cells.Item(5, 4).Value = 3.869799484577842
$ws.Cells.Item(5, 6).Value = 20.16177255877975
$ws.Cells.Item(5, 7).Value = 21.7702218775624
$ws.Cells.Item(5, 8).Value = 12.49306122377543
$ws.Cells.Item(5, 9).Value = 19.39481036204846
$ws.Cells.Item(5, 12).Value = 10.56726432160638
$ws.Cells.Item(5, 14).Value = 16.52898312953872
$ws.Cells.Item(5, 15).Value = 18.08329313974115
$ws.Cells.Item(6, 2).Value = 12.78240257860907
$ws.Cells.Item(6, 3).Value = 11.37377703756553
$ws.Cells.Item(6, 4).Value = 3.866628000594721
$ws.Cells.Item(6, 6).Value = 20.16127000619292
$ws.Cells.Item(6, 7).Value = 21.76786282001633
$ws.Cells.Item(6, 8).Value = 12.49479359086019
$ws.Cells.Item(6, 9).Value = 19.39915595726479
$ws.Cells.Item(6, 12).Value = 10.56629650857453
$ws.Cells.Item(6, 14).Value = 16.53054527615754
$ws.Cells.Item(6, 15).Value = 18.08531774298448
$ws.Cells.Item(7, 2).Value = 12.95574827776673
$ws.Cells.Item(7, 3).Value = 11.40364885426722
$ws.Cells.Item(7, 4).Value = 3.888535071488211
$ws.Cells.Item(7, 6).Value = 20.16532135291771
$ws.Cells.Item(7, 7).Value = 21.78503255632783
$ws.Cells.Item(7, 8).Value = 12.48297168014174
$ws.Cells.Item(7, 9).Value = 19.36937856338499
$ws.Cells.Item(7, 12).Value = 10.57320601960742
$ws.Cells.Item(7, 14).Value = 16.51985387329614
$ws.Cells.Item(7, 15).Value = 18.07168090825772
$ws.Cells.Item(8, 2).Value = 13.69581576888644
$ws.Cells.Item(8, 3).Value = 11.53449645250346
$ws.Cells.Item(8, 4).Value = 3.98247896566957
$ws.Cells.Item(8, 6).Value = 20.19710788715966
$ws.Cells.Item(8, 7).Value = 21.88039013950674
$ws.Cells.Item(8, 8).Value = 12.43582186278695
$ws.Cells.Item(8, 9).Value = 19.24753598783442
$ws.Cells.Item(8, 12).Value = 10.6084100272069
$ws.Cells.Item(8, 14).Value = 16.47642195427327
$ws.Cells.Item(8, 15).Value = 18.02182238400771
$ws.Cells.Item(9, 2).Value = 15.04540001548325
$ws.Cells.Item(9, 3).Value = 11.78653325793697
$ws.Cells.Item(9, 4).Value = 4.156077049121933
$ws.Cells.Item(9, 6).Value = 20.30906843388415
$ws.Cells.Item(9, 7).Value = 22.13665650644751
$ws.Cells.Item(9, 8).Value = 12.36143253409633
$ws.Cells.Item(9, 9).Value = 19.04307679307749
$ws.Cells.Item(9, 12).Value = 10.69403393976173
$ws.Cells.Item(9, 14).Value = 16.4046770082027
$ws.Cells.Item(9, 15).Value = 17.96130460460294
$ws.Cells.Item(10, 2).Value = 15.96608505375053
$ws.Cells.Item(10, 3).Value = 11.96798731675119
$ws.Cells.Item(10, 4).Value = 4.276466600038389
$ws.Cells.Item(10, 6).Value = 20.42111083086978
$ws.Cells.Item(10, 7).Value = 22.36580195270198
$ws.Cells.Item(10, 8).Value = 12.31782696126215
$ws.Cells.Item(10, 9).Value = 18.91398582315201
$ws.Cells.Item(10, 12).Value = 10.76669623381508
$ws.Cells.Item(10, 14).Value = 16.36011268739894
$ws.Cells.Item(10, 15).Value = 17.93967803596224
$ws.Cells.Item(11, 2).Value = 16.36807398589573
$ws.Cells.Item(11, 3).Value = 12.04953205768526
$ws.Cells.Item(11, 4).Value = 4.329551915850995
$ws.Cells.Item(11, 6).Value = 20.47841422168748
$ws.Cells.Item(11, 7).Value = 22.47855542360213
$ws.Cells.Item(11, 8).Value = 12.3003945917102
$ws.Cells.Item(11, 9).Value = 18.85987010651861
$ws.Cells.Item(11, 12).Value = 10.80178358142286
$ws.Cells.Item(11, 14).Value = 16.34160038353485
$ws.Cells.Item(11, 15).Value = 17.93481546507856
$ws.Cells.Item(12, 2).Value = 16.51777639957509
$ws.Cells.Item(12, 3).Value = 12.08025095846045
$ws.Cells.Item(12, 4).Value = 4.349402283414849
$ws.Cells.Item(12, 6).Value = 20.501010483682
$ws.Cells.Item(12, 7).Value = 22.52244109718073
$ws.Cells.Item(12, 8).Value = 12.29413951551266
$ws.Cells.Item(12, 9).Value = 18.84004268802413
$ws.Cells.Item(12, 12).Value = 10.81535451527712
$ws.Cells.Item(12, 14).Value = 16.33484276229051
$ws.Cells.Item(12, 15).Value = 17.93369018636072
$ws.Cells.Item(13, 2).Value = 16.48564892790218
$ws.Cells.Item(13, 3).Value = 12.07364248449896
$ws.Cells.Item(13, 4).Value = 4.345138510272657
$ws.Cells.Item(13, 6).Value = 20.49610433599067
$ws.Cells.Item(13, 7).Value = 22.51293728148462
$ws.Cells.Item(13, 8).Value = 12.29547125428527
$ws.Cells.Item(13, 9).Value = 18.84428326245878
$ws.Cells.Item(13, 12).Value = 10.81241927101042
$ws.Cells.Item(13, 14).Value = 16.33628691151719
$ws.Cells.Item(13, 15).Value = 17.93390068355458
$ws.Cells.Item(14, 2).Value = 16.38044116573822
$ws.Cells.Item(14, 3).Value = 12.0520626270703
$ws.Cells.Item(14, 4).Value = 4.331190122016639
$ws.Cells.Item(14, 6).Value = 20.48025534536179
$ws.Cells.Item(14, 7).Value = 22.48214232358023
$ws.Cells.Item(14, 8).Value = 12.29987304218734
$ws.Cells.Item(14, 9).Value = 18.85822555119129
$ws.Cells.Item(14, 12).Value = 10.80289442476169
$ws.Cells.Item(14, 14).Value = 16.34103937079363
$ws.Cells.Item(14, 15).Value = 17.93470853491763
$ws.Cells.Item(15, 2).Value = 16.31566718044349
$ws.Cells.Item(15, 3).Value = 12.03882297592797
$ws.Cells.Item(15, 4).Value = 4.322613235710087
$ws.Cells.Item(15, 6).Value = 20.47066371038763
$ws.Cells.Item(15, 7).Value = 22.46343316882865
$ws.Cells.Item(15, 8).Value = 12.30261436127708
$ws.Cells.Item(15, 9).Value = 18.86685228195772
$ws.Cells.Item(15, 12).Value = 10.79709693071321
$ws.Cells.Item(15, 14).Value = 16.3439832678098
$ws.Cells.Item(15, 15).Value = 17.93529663016002
$ws.Cells.Item(16, 2).Value = 15.93946636391652
$ws.Cells.Item(16, 3).Value = 11.96263663691474
$ws.Cells.Item(16, 4).Value = 4.272962643538342
$ws.Cells.Item(16, 6).Value = 20.41749221174659
$ws.Cells.Item(16, 7).Value = 22.35860156898627
$ws.Cells.Item(16, 8).Value = 12.3190146256889
$ws.Cells.Item(16, 9).Value = 18.91761534550605
$ws.Cells.Item(16, 12).Value = 10.7644434083079
$ws.Cells.Item(16, 14).Value = 16.36135788470047
$ws.Cells.Item(16, 15).Value = 17.94009599877065
$ws.Cells.Item(17, 2).Value = 15.70429051365055
$ws.Cells.Item(17, 3).Value = 11.91563086559552
$ws.Cells.Item(17, 4).Value = 4.242065701199803
$ws.Cells.Item(17, 6).Value = 20.38648609110884
$ws.Cells.Item(17, 7).Value = 22.2964468180001
$ws.Cells.Item(17, 8).Value = 12.32969174946489
$ws.Cells.Item(17, 9).Value = 18.94993882743104
$ws.Cells.Item(17, 12).Value = 10.74492650213195
$ws.Cells.Item(17, 14).Value = 16.37246710617978
$ws.Cells.Item(17, 15).Value = 17.94431518257172
$ws.Cells.Item(18, 2).Value = 15.56744397358824
$ws.Cells.Item(18, 3).Value = 11.88850127349576
$ws.Cells.Item(18, 4).Value = 4.224136990122237
$ws.Cells.Item(18, 6).Value = 20.36924940052589
$ws.Cells.Item(18, 7).Value = 22.26150078068444
$ws.Cells.Item(18, 8).Value = 12.33605920097417
$ws.Cells.Item(18, 9).Value = 18.96896406332104
$ws.Cells.Item(18, 12).Value = 10.73389276072755
$ws.Cells.Item(18, 14).Value = 16.37902254918048
$ws.Cells.Item(18, 15).Value = 17.94721025516597
$ws.Cells.Item(19, 2).Value = 15.52084212193784
$ws.Cells.Item(19, 3).Value = 11.87930017727068
$ws.Cells.Item(19, 4).Value = 4.21803989039589
$ws.Cells.Item(19, 6).Value = 20.3635163317164
$ws.Cells.Item(19, 7).Value = 22.24980774058963
$ws.Cells.Item(19, 8).Value = 12.33825395394789
$ws.Cells.Item(19, 9).Value = 18.97548008555639
$ws.Cells.Item(19, 12).Value = 10.73019012184611
$ws.Cells.Item(19, 14).Value = 16.38127058714388
$ws.Cells.Item(19, 15).Value = 17.948270882932
$ws.Cells.Item(20, 2).Value = 15.72948964933009
$ws.Cells.Item(20, 3).Value = 11.92064446840698
$ws.Cells.Item(20, 4).Value = 4.245371121473146
$ws.Cells.Item(20, 6).Value = 20.38972503453247
$ws.Cells.Item(20, 7).Value = 22.30298036505232
$ws.Cells.Item(20, 8).Value = 12.32853173143255
$ws.Cells.Item(20, 9).Value = 18.94645304221703
$ws.Cells.Item(20, 12).Value = 10.74698431072956
$ws.Cells.Item(20, 14).Value = 16.37126736381922
$ws.Cells.Item(20, 15).Value = 17.94381756895995
$ws.Cells.Item(21, 2).Value = 16.41141241585608
$ws.Cells.Item(21, 3).Value = 12.05840563595042
$ws.Cells.Item(21, 4).Value = 4.335294016208902
$ws.Cells.Item(21, 6).Value = 20.48488636271712
$ws.Cells.Item(21, 7).Value = 22.49115559762133
$ws.Cells.Item(21, 8).Value = 12.29857073143295
$ws.Cells.Item(21, 9).Value = 18.85411229554052
$ws.Cells.Item(21, 12).Value = 10.80568445972161
$ws.Cells.Item(21, 14).Value = 16.33963660740077
$ws.Cells.Item(21, 15).Value = 17.93445181310526
$ws.Cells.Item(22, 2).Value = 16.84235111698269
$ws.Cells.Item(22, 3).Value = 12.14749858264434
$ws.Cells.Item(22, 4).Value = 4.392591346581975
$ws.Cells.Item(22, 6).Value = 20.55229882048068
$ws.Cells.Item(22, 7).Value = 22.62104712955562
$ws.Cells.Item(22, 8).Value = 12.28100749733801
$ws.Cells.Item(22, 9).Value = 18.7976396793336
$ws.Cells.Item(22, 12).Value = 10.84570034467827
$ws.Cells.Item(22, 14).Value = 16.32043610260769
$ws.Cells.Item(22, 15).Value = 17.93250448306137
$ws.Cells.Item(23, 2).Value = 16.61372910179739
$ws.Cells.Item(23, 3).Value = 12.10003958015338
$ws.Cells.Item(23, 4).Value = 4.362148623263792
$ws.Cells.Item(23, 6).Value = 20.51584708227114
$ws.Cells.Item(23, 7).Value = 22.55110208621449
$ws.Cells.Item(23, 8).Value = 12.29019653080474
$ws.Cells.Item(23, 9).Value = 18.82742461511357
$ws.Cells.Item(23, 12).Value = 10.82419478336882
$ws.Cells.Item(23, 14).Value = 16.33054925799166
$ws.Cells.Item(23, 15).Value = 17.93316184428197
$ws.Cells.Item(24, 2).Value = 15.71810222300689
$ws.Cells.Item(24, 3).Value = 11.91837814545452
$ws.Cells.Item(24, 4).Value = 4.243877256140233
$ws.Cells.Item(24, 6).Value = 20.3882588720255
$ws.Cells.Item(24, 7).Value = 22.30002409256343
$ws.Cells.Item(24, 8).Value = 12.32905546220151
$ws.Cells.Item(24, 9).Value = 18.94802758897278
$ws.Cells.Item(24, 12).Value = 10.74605339302563
$ws.Cells.Item(24, 14).Value = 16.37180924210827
$ws.Cells.Item(24, 15).Value = 17.94404107795104
$ws.Cells.Item(25, 2).Value = 14.69218468101464
$ws.Cells.Item(25, 3).Value = 11.71893595892688
$ws.Cells.Item(25, 4).Value = 4.110323973085261
$ws.Cells.Item(25, 6).Value = 20.27350890846638
$ws.Cells.Item(25, 7).Value = 22.06004016982896
$ws.Cells.Item(25, 8).Value = 12.37961872671726
$ws.Cells.Item(25, 9).Value = 19.09468660812742
$ws.Cells.Item(25, 12).Value = 10.69403393976173
$ws.Cells.Item(25, 14).Value = 16.42265258448151
$ws.Cells.Item(25, 15).Value = 17.97367326792174
